# chore: update Sheets via scheduled runner
# Refreshes market-price derived columns (H-N) for several Leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1344.6154
$ws.Range("I18").Value = 1172.5
$ws.Range("J18").Value = 1620
$ws.Range("K18").Value = 1172.5
$ws.Range("L18").Value = 1620
$ws.Range("M18").Value = -888.5
$ws.Range("N18").Value = -2188

$ws.Range("H34").Value = 13024
$ws.Range("I34").Value = 1628.8
$ws.Range("J34").Value = 70000
$ws.Range("K34").Value = 1628.8
$ws.Range("L34").Value = 70000
$ws.Range("M34").Value = -1425.8
$ws.Range("N34").Value = -70406

$ws.Range("H36").Value = 13024
$ws.Range("I36").Value = 1628.8
$ws.Range("J36").Value = 70000
$ws.Range("K36").Value = 1628.8
$ws.Range("L36").Value = 70000
$ws.Range("M36").Value = -913.8
$ws.Range("N36").Value = -71430

$ws.Range("H137").Value = 39052.445
$ws.Range("I137").Value = 1206.875
$ws.Range("J137").Value = 94100.55
$ws.Range("K137").Value = 3620.625
$ws.Range("L137").Value = 282301.65
$ws.Range("M137").Value = -1070.625
$ws.Range("N137").Value = -287401.65

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4818.093
$ws.Range("I32").Value = 4664.3022
$ws.Range("K32").Value = 4664.3022
$ws.Range("M32").Value = -4377.3022

$ws.Range("H74").Value = 5328.25
$ws.Range("I74").Value = 638.8889
$ws.Range("J74").Value = 8141.8667
$ws.Range("K74").Value = 638.8889
$ws.Range("L74").Value = 8141.8667
$ws.Range("M74").Value = 235.1111
$ws.Range("N74").Value = -9889.866699999999

$ws.Range("H77").Value = 5328.25
$ws.Range("I77").Value = 638.8889
$ws.Range("J77").Value = 8141.8667
$ws.Range("K77").Value = 3194.4445
$ws.Range("L77").Value = 40709.3335
$ws.Range("M77").Value = 1173.5555
$ws.Range("N77").Value = -49445.3335

$ws.Range("H110").Value = 846.9091
$ws.Range("I110").Value = 741.6
$ws.Range("J110").Value = 1900
$ws.Range("K110").Value = 741.6
$ws.Range("L110").Value = 1900
$ws.Range("M110").Value = 1303.4
$ws.Range("N110").Value = -5990

$ws.Range("H133").Value = 46235
$ws.Range("J133").Value = 46235
$ws.Range("L133").Value = 46235
$ws.Range("N133").Value = -51295

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 79780
$ws.Range("J59").Value = 79780
$ws.Range("L59").Value = 79780
$ws.Range("N59").Value = -81474

$ws.Range("H122").Value = 14800
$ws.Range("J122").Value = 14800
$ws.Range("L122").Value = 14800
$ws.Range("N122").Value = -24600

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 23126.584
$ws.Range("I41").Value = 6029.5
$ws.Range("J41").Value = 26546
$ws.Range("K41").Value = 6029.5
$ws.Range("L41").Value = 26546
$ws.Range("M41").Value = -5601.5
$ws.Range("N41").Value = -27402

$ws.Range("H50").Value = 32478.857
$ws.Range("I50").Value = 12000
$ws.Range("J50").Value = 35892
$ws.Range("K50").Value = 12000
$ws.Range("L50").Value = 35892
$ws.Range("M50").Value = -11375
$ws.Range("N50").Value = -37142

$ws.Range("H51").Value = 24283.77
$ws.Range("I51").Value = 12000
$ws.Range("J51").Value = 25307.416
$ws.Range("K51").Value = 12000
$ws.Range("L51").Value = 25307.416
$ws.Range("M51").Value = -11264
$ws.Range("N51").Value = -26779.416

$ws.Range("H59").Value = 35267.5
$ws.Range("J59").Value = 35267.5
$ws.Range("L59").Value = 35267.5
$ws.Range("N59").Value = -37557.5

$ws.Range("H60").Value = 35501.875
$ws.Range("I60").Value = 21750
$ws.Range("J60").Value = 40085.832
$ws.Range("K60").Value = 21750
$ws.Range("L60").Value = 40085.832
$ws.Range("M60").Value = -21239
$ws.Range("N60").Value = -41107.832

$ws.Range("H61").Value = 24283.77
$ws.Range("I61").Value = 12000
$ws.Range("J61").Value = 25307.416
$ws.Range("K61").Value = 12000
$ws.Range("L61").Value = 25307.416
$ws.Range("M61").Value = -11652
$ws.Range("N61").Value = -26003.416

$ws.Range("H68").Value = 44516.668
$ws.Range("J68").Value = 47596.152
$ws.Range("L68").Value = 47596.152
$ws.Range("N68").Value = -49094.152

$ws.Range("H71").Value = 44516.668
$ws.Range("J71").Value = 47596.152
$ws.Range("L71").Value = 142788.456
$ws.Range("N71").Value = -150276.456

$ws.Range("H74").Value = 50314
$ws.Range("J74").Value = 50314
$ws.Range("L74").Value = 50314
$ws.Range("N74").Value = -52062

$ws.Range("H77").Value = 50314
$ws.Range("J77").Value = 50314
$ws.Range("L77").Value = 150942
$ws.Range("N77").Value = -159678

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3565.7673
$ws.Range("I68").Value = 570
$ws.Range("J68").Value = 3959.9473
$ws.Range("K68").Value = 1710
$ws.Range("L68").Value = 11879.8419
$ws.Range("M68").Value = -899
$ws.Range("N68").Value = -13501.8419

$ws.Range("H71").Value = 3565.7673
$ws.Range("I71").Value = 570
$ws.Range("J71").Value = 3959.9473
$ws.Range("K71").Value = 5130
$ws.Range("L71").Value = 35639.5257
$ws.Range("M71").Value = -1074
$ws.Range("N71").Value = -43751.5257

$ws.Range("H96").Value = 3517.7778
$ws.Range("I96").Value = 2800
$ws.Range("J96").Value = 3876.6667
$ws.Range("K96").Value = 8400
$ws.Range("L96").Value = 11630.0001
$ws.Range("M96").Value = -6341
$ws.Range("N96").Value = -15748.0001

$ws.Range("H107").Value = 4168.625
$ws.Range("J107").Value = 1271.72
$ws.Range("L107").Value = 3815.16
$ws.Range("N107").Value = -7655.16

$ws.Range("H137").Value = 4202.636
$ws.Range("I137").Value = 648
$ws.Range("J137").Value = 4558.1
$ws.Range("K137").Value = 1944
$ws.Range("L137").Value = 13674.3
$ws.Range("M137").Value = 3156
$ws.Range("N137").Value = -23874.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1394.2
$ws.Range("I122").Value = 1201.1666
$ws.Range("J122").Value = 2166.3333
$ws.Range("K122").Value = 3603.4998
$ws.Range("L122").Value = 6498.999899999999
$ws.Range("M122").Value = -1153.4998
$ws.Range("N122").Value = -11398.9999

$ws.Range("H126").Value = 1006
$ws.Range("I126").Value = 1006
$ws.Range("K126").Value = 3018
$ws.Range("M126").Value = -548

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1668033.9
$ws.Range("I126").Value = 2001040.6
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 6003121.800000001
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -6000651.800000001
$ws.Range("N126").Value = -13940
